$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sheet view: topLeftCell and selection
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("E20").Select()

# Columns N..Z (14 columns) map to formulas that concatenate column A (absolute
# column reference) of the corresponding source row with the value in the same
# column of that source row, followed by a trailing comma.
$cols = @("N","O","P","Q","R","S","T","U","V","W","X","Y","Z")

for ($k = 0; $k -le 14; $k++) {
    $targetRow = 19 + $k
    $sourceRow = 1 + $k
    foreach ($col in $cols) {
        $formula = "=`$A" + $sourceRow + "&"":'""&" + $col + $sourceRow + "&""',"""
        $ws.Range($col + $targetRow).Formula = $formula
    }
}
